# "With continuous analysis and graphs"
#
# The e1/e2 outcome columns ("82.5 (20-200)", "45 (10-200)", ...) are split
# into a plain mean/value column plus a separate "sd/mean" (sd or range)
# column, so two new columns are inserted: one right after the e1 column
# (G) and one right after the e2 column (which, after the first insertion,
# sits at column I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new columns - this shifts the existing e2/p-value columns
# two places to the right (H,I -> I,K) and leaves two blank columns (H,J)
# for the new "sd/mean" data.
$ws.Columns("H").Insert()
$ws.Columns("J").Insert()

# New header row
$ws.Range("H1").Value = "sd/mean"
$ws.Range("J1").Value = "sd/mean"

# Row 2 - Crippa et.al., 2008: "82.5 (20-200)" / "45 (10-200)"
# "82.5" still looks like a plain number, so force the cell to Text first or
# Excel's automatic type detection would silently convert it back to a number.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "82.5"
$ws.Range("H2").Value = "20-200"
$ws.Range("I2").Value = 45
$ws.Range("J2").Value = "10-200"

# Row 3 - Yamao, K., et al 2011: "60.1 (38.0)" / "90.0 (45.8)"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "60.1"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "38.0"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "90.0"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "45.8"

# Row 4 - Gil, E., et al., 2012: "56 (21)" / "52 (31)"
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 21
$ws.Range("I4").Value = 52
$ws.Range("J4").Value = 31

# Row 5 - Jang, K. T., et al (2015): values were already bare numbers (94 / 54),
# nothing to split - the sd/mean cells stay empty (fully cleared, not just blank).
$ws.Range("G5").Value = 94
$ws.Range("I5").Value = 54
$ws.Range("H5").Clear()
$ws.Range("J5").Clear()

# Row 6 - Keane 2018: "100 (45-131) " / "B= 52 (30-85) "
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = "45–131"
$ws.Range("I6").Value = 52
$ws.Range("J6").Value = "30–85"
